$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change the header text in B1 from "Invoice" to "Lease Agreement"
$ws.Range("B1").Value = "Lease Agreement"

# Update the active cell selection to B2
$ws.Range("B2").Select()

# Widen column B (4 characters wider than the other columns) to fit the new header text
$ws.Columns("B").ColumnWidth = $ws.Columns("A").ColumnWidth + 4

# Increase row 1 height to accommodate the taller header
$ws.Rows(1).RowHeight = 18.75
